# Update TPM-derived values in the LR-pairs sheet (Mstn-Acvr2b).
#
# The underlying TPM data changed, which updated the Receptor average/total
# expression values in row 2 (M2, N2). That cascades into the derived
# specificity and edge-weight columns (O-T) for every row, because those
# columns are computed relative to the sum of the Receptor expression
# columns across all rows:
#   O (Receptor derived specificity of avg expr)  = M / sum(M)
#   P = O
#   Q (Edge average expression weight)            = Ligand avg * Receptor avg (G * M)
#   R (Edge total expression weight)               = Ligand total * Receptor total (H * N)
#   S (Edge average expression derived specificity)= Q / sum(Q)
#   T (Edge total expression derived specificity)  = R / sum(R)
#
# The values below are the recalculated results (matching the committed
# workbook) so the numbers land exactly on the expected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> ECs)
$ws.Range("M2").Value2 = 1.485259333333333
$ws.Range("N2").Value2 = 4.455778
$ws.Range("O2").Value2 = 0.3057455162066235
$ws.Range("P2").Value2 = 0.3057455162066235
$ws.Range("Q2").Value2 = 0.0627160655228889
$ws.Range("R2").Value2 = 0.5644445897060001
$ws.Range("S2").Value2 = 0.3057455162066235
$ws.Range("T2").Value2 = 0.3057455162066235

# Row 3 (MuSCs -> FAPs) - M3/N3 unchanged, but their relative share shifted
$ws.Range("O3").Value2 = 0.2805555239151429
$ws.Range("P3").Value2 = 0.2805555239151429
$ws.Range("S3").Value2 = 0.2805555239151429
$ws.Range("T3").Value2 = 0.2805555239151429

# Row 4 (MuSCs -> MuSCs) - M4/N4 unchanged, but their relative share shifted
$ws.Range("O4").Value2 = 0.4136989598782336
$ws.Range("P4").Value2 = 0.4136989598782336
$ws.Range("S4").Value2 = 0.4136989598782336
$ws.Range("T4").Value2 = 0.4136989598782336
